$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 531.8889
$ws.Range("I4").Value = 398.6
$ws.Range("K4").Value = 398.6
$ws.Range("M4").Value = -284.6
$ws.Range("H6").Value = 10250
$ws.Range("I6").Value = 10250
$ws.Range("K6").Value = 30750
$ws.Range("M6").Value = -30638
$ws.Range("H11").Value = 865.4
$ws.Range("I11").Value = 865.4
$ws.Range("K11").Value = 865.4
$ws.Range("M11").Value = -725.4
$ws.Range("H28").Value = 562
$ws.Range("I28").Value = 418
$ws.Range("J28").Value = 994
$ws.Range("K28").Value = 418
$ws.Range("L28").Value = 994
$ws.Range("M28").Value = 67
$ws.Range("N28").Value = -1964
$ws.Range("H38").Value = 3441.5557
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30744
$ws.Range("H39").Value = 404.7353
$ws.Range("I39").Value = 185.54546
$ws.Range("J39").Value = 509.56522
$ws.Range("K39").Value = 556.6363799999999
$ws.Range("L39").Value = 1528.69566
$ws.Range("M39").Value = -260.6363799999999
$ws.Range("N39").Value = -2120.69566
$ws.Range("H43").Value = 9224.5
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 11632.667
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 11632.667
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -11770.667
$ws.Range("H86").Value = 5212.1787
$ws.Range("I86").Value = 5404.7896
$ws.Range("J86").Value = 4805.5557
$ws.Range("K86").Value = 5404.7896
$ws.Range("L86").Value = 4805.5557
$ws.Range("M86").Value = -4281.7896
$ws.Range("N86").Value = -7051.5557
$ws.Range("H87").Value = 120000
$ws.Range("J87").Value = 120000
$ws.Range("L87").Value = 120000
$ws.Range("N87").Value = -122496
$ws.Range("H88").Value = 519843.38
$ws.Range("I88").Value = 51300.5
$ws.Range("J88").Value = 676024.3
$ws.Range("K88").Value = 51300.5
$ws.Range("L88").Value = 676024.3
$ws.Range("M88").Value = -50894.5
$ws.Range("N88").Value = -676836.3
$ws.Range("H89").Value = 5212.1787
$ws.Range("I89").Value = 5404.7896
$ws.Range("J89").Value = 4805.5557
$ws.Range("K89").Value = 27023.948
$ws.Range("L89").Value = 24027.7785
$ws.Range("M89").Value = -21407.948
$ws.Range("N89").Value = -35259.7785
$ws.Range("H90").Value = 120000
$ws.Range("J90").Value = 120000
$ws.Range("L90").Value = 360000
$ws.Range("N90").Value = -372480
$ws.Range("H91").Value = 519843.38
$ws.Range("I91").Value = 51300.5
$ws.Range("J91").Value = 676024.3
$ws.Range("K91").Value = 51300.5
$ws.Range("L91").Value = 676024.3
$ws.Range("M91").Value = -49896.5
$ws.Range("N91").Value = -678832.3
$ws.Range("H98").Value = 90911400
$ws.Range("I98").Value = 100002260
$ws.Range("K98").Value = 100002260
$ws.Range("M98").Value = -100000762
$ws.Range("H107").Value = 1226.4445
$ws.Range("I107").Value = 685.2
$ws.Range("J107").Value = 3932.6667
$ws.Range("K107").Value = 685.2
$ws.Range("L107").Value = 3932.6667
$ws.Range("M107").Value = 1234.8
$ws.Range("N107").Value = -7772.6667
$ws.Range("H113").Value = 83336340
$ws.Range("I113").Value = 33336000
$ws.Range("J113").Value = 133336664
$ws.Range("K113").Value = 33336000
$ws.Range("L113").Value = 133336664
$ws.Range("M113").Value = -33332746
$ws.Range("N113").Value = -133343172
$ws.Range("H115").Value = 121401.5
$ws.Range("I115").Value = 121401.5
$ws.Range("K115").Value = 364204.5
$ws.Range("M115").Value = -362637.5
$ws.Range("H116").Value = 5738.4443
$ws.Range("I116").Value = 5395.5
$ws.Range("K116").Value = 5395.5
$ws.Range("M116").Value = -1953.5
$ws.Range("H122").Value = 90911400
$ws.Range("I122").Value = 100002260
$ws.Range("K122").Value = 300006780
$ws.Range("M122").Value = -300004330
$ws.Range("H132").Value = 2409.35
$ws.Range("I132").Value = 2349
$ws.Range("J132").Value = 2952.5
$ws.Range("K132").Value = 7047
$ws.Range("L132").Value = 8857.5
$ws.Range("M132").Value = -4517
$ws.Range("N132").Value = -13917.5
$ws.Range("H135").Value = 10979.923
$ws.Range("J135").Value = 22574.4
$ws.Range("L135").Value = 203169.6
$ws.Range("N135").Value = -208239.6
$ws.Range("H137").Value = 3626.7812
$ws.Range("I137").Value = 2440.7
$ws.Range("K137").Value = 7322.099999999999
$ws.Range("M137").Value = -4772.099999999999
$ws.Range("H138").Value = 1013183.1
$ws.Range("I138").Value = 1439.2
$ws.Range("J138").Value = 1193851.6
$ws.Range("K138").Value = 4317.6
$ws.Range("L138").Value = 3581554.8
$ws.Range("M138").Value = 822.3999999999996
$ws.Range("N138").Value = -3591834.8
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1033.5
$ws.Range("I25").Value = 1040.2
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 1040.2
$ws.Range("L25").Value = 1000
$ws.Range("M25").Value = -638.2
$ws.Range("N25").Value = -1804
$ws.Range("H32").Value = 8475998
$ws.Range("I32").Value = 8622134
$ws.Range("J32").Value = 95
$ws.Range("K32").Value = 8622134
$ws.Range("L32").Value = 95
$ws.Range("M32").Value = -8621847
$ws.Range("N32").Value = -669
$ws.Range("H45").Value = 2293.0527
$ws.Range("I45").Value = 2218.25
$ws.Range("J45").Value = 2313
$ws.Range("K45").Value = 2218.25
$ws.Range("L45").Value = 2313
$ws.Range("M45").Value = -1841.25
$ws.Range("N45").Value = -3067
$ws.Range("H61").Value = 41756970
$ws.Range("I61").Value = 125003660
$ws.Range("K61").Value = 125003660
$ws.Range("M61").Value = -125003448
$ws.Range("H74").Value = 5324713.5
$ws.Range("I74").Value = 8065986
$ws.Range("J74").Value = 13498.25
$ws.Range("K74").Value = 8065986
$ws.Range("L74").Value = 13498.25
$ws.Range("M74").Value = -8065112
$ws.Range("N74").Value = -15246.25
$ws.Range("H76").Value = 55396
$ws.Range("J76").Value = 55396
$ws.Range("L76").Value = 55396
$ws.Range("N76").Value = -56072
$ws.Range("H77").Value = 5324713.5
$ws.Range("I77").Value = 8065986
$ws.Range("J77").Value = 13498.25
$ws.Range("K77").Value = 40329930
$ws.Range("L77").Value = 67491.25
$ws.Range("M77").Value = -40325562
$ws.Range("N77").Value = -76227.25
$ws.Range("H79").Value = 55396
$ws.Range("J79").Value = 55396
$ws.Range("L79").Value = 55396
$ws.Range("N79").Value = -57736
$ws.Range("H82").Value = 39534.5
$ws.Range("J82").Value = 39534.5
$ws.Range("L82").Value = 39534.5
$ws.Range("N82").Value = -40256.5
$ws.Range("H85").Value = 39534.5
$ws.Range("J85").Value = 39534.5
$ws.Range("L85").Value = 39534.5
$ws.Range("N85").Value = -42030.5
$ws.Range("H88").Value = 1864.909
$ws.Range("J88").Value = 1830.8572
$ws.Range("L88").Value = 1830.8572
$ws.Range("N88").Value = -2642.8572
$ws.Range("H91").Value = 1864.909
$ws.Range("J91").Value = 1830.8572
$ws.Range("L91").Value = 1830.8572
$ws.Range("N91").Value = -4638.8572
$ws.Range("H97").Value = 1287.3846
$ws.Range("I97").Value = 875.3
$ws.Range("K97").Value = 875.3
$ws.Range("M97").Value = -379.3
$ws.Range("H102").Value = 19662.5
$ws.Range("I102").Value = 19662.5
$ws.Range("K102").Value = 19662.5
$ws.Range("M102").Value = -18040.5
$ws.Range("H110").Value = 2079.5334
$ws.Range("I110").Value = 1956.381
$ws.Range("J110").Value = 2366.889
$ws.Range("K110").Value = 1956.381
$ws.Range("L110").Value = 2366.889
$ws.Range("M110").Value = 88.61899999999991
$ws.Range("N110").Value = -6456.889
$ws.Range("H122").Value = 4445.778
$ws.Range("I122").Value = 3337.3333
$ws.Range("K122").Value = 10011.9999
$ws.Range("M122").Value = -7561.999899999999
$ws.Range("H132").Value = 4715.684
$ws.Range("I132").Value = 2090.9312
$ws.Range("J132").Value = 13173.223
$ws.Range("K132").Value = 6272.7936
$ws.Range("L132").Value = 39519.669
$ws.Range("M132").Value = -3742.7936
$ws.Range("N132").Value = -44579.669
$ws.Range("H136").Value = 41756970
$ws.Range("I136").Value = 125003660
$ws.Range("K136").Value = 375010980
$ws.Range("M136").Value = -375008430
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1471
$ws.Range("I86").Value = 1416.3
$ws.Range("J86").Value = 1580.4
$ws.Range("K86").Value = 1416.3
$ws.Range("L86").Value = 1580.4
$ws.Range("M86").Value = -293.3
$ws.Range("N86").Value = -3826.4
$ws.Range("H89").Value = 1471
$ws.Range("I89").Value = 1416.3
$ws.Range("J89").Value = 1580.4
$ws.Range("K89").Value = 7081.5
$ws.Range("L89").Value = 7902
$ws.Range("M89").Value = -1465.5
$ws.Range("N89").Value = -19134
$ws.Range("H94").Value = 1921.5264
$ws.Range("I94").Value = 1292.4166
$ws.Range("K94").Value = 1292.4166
$ws.Range("M94").Value = -841.4166
$ws.Range("H99").Value = 2252.4666
$ws.Range("I99").Value = 1439
$ws.Range("K99").Value = 1439
$ws.Range("M99").Value = 59
$ws.Range("H105").Value = 2125.7
$ws.Range("I105").Value = 1967
$ws.Range("K105").Value = 1967
$ws.Range("M105").Value = -220
$ws.Range("H107").Value = 1885.1052
$ws.Range("I107").Value = 2154
$ws.Range("J107").Value = 1132.2
$ws.Range("K107").Value = 2154
$ws.Range("L107").Value = 1132.2
$ws.Range("M107").Value = -234
$ws.Range("N107").Value = -4972.2
$ws.Range("H134").Value = 29758.77
$ws.Range("I134").Value = 3888
$ws.Range("K134").Value = 11664
$ws.Range("M134").Value = -9129

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1974.1111
$ws.Range("J16").Value = 2500.25
$ws.Range("L16").Value = 2500.25
$ws.Range("N16").Value = -3074.25
$ws.Range("H22").Value = 468.66666
$ws.Range("I22").Value = 262.66666
$ws.Range("J22").Value = 674.6667
$ws.Range("K22").Value = 262.66666
$ws.Range("L22").Value = 674.6667
$ws.Range("M22").Value = 87.33334000000002
$ws.Range("N22").Value = -1374.6667
$ws.Range("H31").Value = 1062632.1
$ws.Range("I31").Value = 1682.5
$ws.Range("J31").Value = 1460488.2
$ws.Range("K31").Value = 1682.5
$ws.Range("L31").Value = 1460488.2
$ws.Range("M31").Value = -1387.5
$ws.Range("N31").Value = -1461078.2
$ws.Range("H32").Value = 4670
$ws.Range("I32").Value = 4670
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4670
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4354
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 1062632.1
$ws.Range("I34").Value = 1682.5
$ws.Range("J34").Value = 1460488.2
$ws.Range("K34").Value = 1682.5
$ws.Range("L34").Value = 1460488.2
$ws.Range("M34").Value = -1480.5
$ws.Range("N34").Value = -1460892.2
$ws.Range("H58").Value = 1325.8125
$ws.Range("I58").Value = 1479.0834
$ws.Range("J58").Value = 866
$ws.Range("K58").Value = 1479.0834
$ws.Range("L58").Value = 866
$ws.Range("M58").Value = -1276.0834
$ws.Range("N58").Value = -1272
$ws.Range("H86").Value = 3543.875
$ws.Range("I86").Value = 3298.25
$ws.Range("J86").Value = 3789.5
$ws.Range("K86").Value = 3298.25
$ws.Range("L86").Value = 3789.5
$ws.Range("M86").Value = -2175.25
$ws.Range("N86").Value = -6035.5
$ws.Range("H89").Value = 3543.875
$ws.Range("I89").Value = 3298.25
$ws.Range("J89").Value = 3789.5
$ws.Range("K89").Value = 16491.25
$ws.Range("L89").Value = 18947.5
$ws.Range("M89").Value = -10875.25
$ws.Range("N89").Value = -30179.5
$ws.Range("H93").Value = 61350.75
$ws.Range("I93").Value = 45801
$ws.Range("K93").Value = 45801
$ws.Range("M93").Value = -43929
$ws.Range("H99").Value = 3964.875
$ws.Range("I99").Value = 3220.4546
$ws.Range("J99").Value = 5602.6
$ws.Range("K99").Value = 3220.4546
$ws.Range("L99").Value = 5602.6
$ws.Range("M99").Value = -1722.4546
$ws.Range("N99").Value = -8598.6
$ws.Range("H103").Value = 29309.428
$ws.Range("I103").Value = 16532.334
$ws.Range("J103").Value = 52308.2
$ws.Range("K103").Value = 16532.334
$ws.Range("L103").Value = 52308.2
$ws.Range("M103").Value = -15360.334
$ws.Range("N103").Value = -54652.2
$ws.Range("H105").Value = 2012.7222
$ws.Range("I105").Value = 1735.6666
$ws.Range("J105").Value = 2566.8333
$ws.Range("K105").Value = 1735.6666
$ws.Range("L105").Value = 2566.8333
$ws.Range("M105").Value = 11.33339999999998
$ws.Range("N105").Value = -6060.8333
$ws.Range("H107").Value = 852.0833
$ws.Range("I107").Value = 666.6111
$ws.Range("K107").Value = 666.6111
$ws.Range("M107").Value = 1253.3889
$ws.Range("H113").Value = 1974.1111
$ws.Range("J113").Value = 2500.25
$ws.Range("L113").Value = 2500.25
$ws.Range("N113").Value = -6840.25
$ws.Range("H122").Value = 4308.5835
$ws.Range("J122").Value = 5384.3076
$ws.Range("L122").Value = 16152.9228
$ws.Range("N122").Value = -21052.9228
$ws.Range("H126").Value = 3964.875
$ws.Range("I126").Value = 3220.4546
$ws.Range("J126").Value = 5602.6
$ws.Range("K126").Value = 9661.363799999999
$ws.Range("L126").Value = 16807.8
$ws.Range("M126").Value = -7191.363799999999
$ws.Range("N126").Value = -21747.8
$ws.Range("H132").Value = 1991.9412
$ws.Range("I132").Value = 2216.7856
$ws.Range("J132").Value = 942.6667
$ws.Range("K132").Value = 6650.3568
$ws.Range("L132").Value = 2828.0001
$ws.Range("M132").Value = -4120.3568
$ws.Range("N132").Value = -7888.0001
$ws.Range("H134").Value = 305907.5
$ws.Range("I134").Value = 358697.94
$ws.Range("J134").Value = 10281.2
$ws.Range("K134").Value = 1076093.82
$ws.Range("L134").Value = 30843.6
$ws.Range("M134").Value = -1073558.82
$ws.Range("N134").Value = -35913.60000000001
$ws.Range("H136").Value = 1325.8125
$ws.Range("I136").Value = 1479.0834
$ws.Range("J136").Value = 866
$ws.Range("K136").Value = 4437.2502
$ws.Range("L136").Value = 2598
$ws.Range("M136").Value = -1887.2502
$ws.Range("N136").Value = -7698

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 322.5
$ws.Range("I14").Value = 322.5
$ws.Range("K14").Value = 967.5
$ws.Range("M14").Value = -794.5
$ws.Range("H18").Value = 2499.5
$ws.Range("J18").Value = 2499.5
$ws.Range("L18").Value = 7498.5
$ws.Range("N18").Value = -7836.5
$ws.Range("H23").Value = 376.83334
$ws.Range("J23").Value = 337.125
$ws.Range("L23").Value = 1011.375
$ws.Range("N23").Value = -1481.375
$ws.Range("H26").Value = 147
$ws.Range("I26").Value = 138.33333
$ws.Range("J26").Value = 160
$ws.Range("K26").Value = 414.99999
$ws.Range("L26").Value = 480
$ws.Range("M26").Value = -126.99999
$ws.Range("N26").Value = -1056
$ws.Range("H56").Value = 9422.77
$ws.Range("I56").Value = 9422.77
$ws.Range("K56").Value = 9422.77
$ws.Range("M56").Value = -8892.77
$ws.Range("H86").Value = 2646.875
$ws.Range("I86").Value = 725
$ws.Range("K86").Value = 2175
$ws.Range("M86").Value = -989
$ws.Range("H89").Value = 2646.875
$ws.Range("I89").Value = 725
$ws.Range("K89").Value = 6525
$ws.Range("M89").Value = -597
$ws.Range("H92").Value = 5000350
$ws.Range("J92").Value = 700
$ws.Range("L92").Value = 2100
$ws.Range("N92").Value = -4596
$ws.Range("H107").Value = 401.13333
$ws.Range("I107").Value = 356.66666
$ws.Range("J107").Value = 412.25
$ws.Range("K107").Value = 1069.99998
$ws.Range("L107").Value = 1236.75
$ws.Range("M107").Value = 850.0000199999999
$ws.Range("N107").Value = -5076.75
$ws.Range("H131").Value = 1397.1538
$ws.Range("J131").Value = 1900
$ws.Range("L131").Value = 5700
$ws.Range("N131").Value = -15780
$ws.Range("H132").Value = 2187.7
$ws.Range("I132").Value = 2137.7273
$ws.Range("J132").Value = 2248.7778
$ws.Range("K132").Value = 19239.5457
$ws.Range("L132").Value = 20239.0002
$ws.Range("M132").Value = -16709.5457
$ws.Range("N132").Value = -25299.0002
$ws.Range("H137").Value = 8479
$ws.Range("I137").Value = 4644.25
$ws.Range("J137").Value = 10396.375
$ws.Range("K137").Value = 13932.75
$ws.Range("L137").Value = 31189.125
$ws.Range("M137").Value = -8832.75
$ws.Range("N137").Value = -41389.125
$ws.Range("H138").Value = 1564.6666
$ws.Range("I138").Value = 1564.6666
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4693.9998
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 446.0002000000004
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 2236.8948
$ws.Range("I139").Value = 1785.8572
$ws.Range("K139").Value = 5357.571599999999
$ws.Range("M139").Value = -217.5715999999993
$ws.Range("H141").Value = 12278.608
$ws.Range("I141").Value = 9580.9
$ws.Range("K141").Value = 28742.7
$ws.Range("M141").Value = -23562.7

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6766.8823
$ws.Range("I70").Value = 4375.75
$ws.Range("K70").Value = 4375.75
$ws.Range("M70").Value = -4105.75
$ws.Range("H73").Value = 6766.8823
$ws.Range("I73").Value = 4375.75
$ws.Range("K73").Value = 4375.75
$ws.Range("M73").Value = -3439.75
$ws.Range("H80").Value = 7224.6113
$ws.Range("I80").Value = 2862.2
$ws.Range("J80").Value = 12677.625
$ws.Range("K80").Value = 2862.2
$ws.Range("L80").Value = 12677.625
$ws.Range("M80").Value = -1864.2
$ws.Range("N80").Value = -14673.625
$ws.Range("H83").Value = 7224.6113
$ws.Range("I83").Value = 2862.2
$ws.Range("J83").Value = 12677.625
$ws.Range("K83").Value = 14311
$ws.Range("L83").Value = 63388.125
$ws.Range("M83").Value = -9319
$ws.Range("N83").Value = -73372.125
$ws.Range("H97").Value = 1076.4762
$ws.Range("I97").Value = 939.6875
$ws.Range("J97").Value = 1514.2
$ws.Range("K97").Value = 939.6875
$ws.Range("L97").Value = 1514.2
$ws.Range("M97").Value = -443.6875
$ws.Range("N97").Value = -2506.2
$ws.Range("H102").Value = 3215.6956
$ws.Range("I102").Value = 2683.3125
$ws.Range("J102").Value = 4432.5713
$ws.Range("K102").Value = 2683.3125
$ws.Range("L102").Value = 4432.5713
$ws.Range("M102").Value = -1061.3125
$ws.Range("N102").Value = -7676.5713
$ws.Range("H113").Value = 4352.4
$ws.Range("I113").Value = 2749.5
$ws.Range("J113").Value = 4753.125
$ws.Range("K113").Value = 2749.5
$ws.Range("L113").Value = 4753.125
$ws.Range("M113").Value = -579.5
$ws.Range("N113").Value = -9093.125
$ws.Range("H126").Value = 4927.7144
$ws.Range("I126").Value = 5070.857
$ws.Range("J126").Value = 4784.5713
$ws.Range("K126").Value = 15212.571
$ws.Range("L126").Value = 14353.7139
$ws.Range("M126").Value = -12742.571
$ws.Range("N126").Value = -19293.7139
$ws.Range("H132").Value = 26317998
$ws.Range("I132").Value = 27779848
$ws.Range("K132").Value = 83339544
$ws.Range("M132").Value = -83337014
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 10125
$ws.Range("I3").Value = 250
$ws.Range("K3").Value = 250
$ws.Range("M3").Value = -138
$ws.Range("H7").Value = 63311.35
$ws.Range("I7").Value = 3564.4443
$ws.Range("K7").Value = 3564.4443
$ws.Range("M7").Value = -3452.4443
$ws.Range("H15").Value = 10125
$ws.Range("I15").Value = 250
$ws.Range("K15").Value = 250
$ws.Range("M15").Value = -80
$ws.Range("H16").Value = 2666.7222
$ws.Range("I16").Value = 2439.2727
$ws.Range("J16").Value = 3024.1428
$ws.Range("K16").Value = 2439.2727
$ws.Range("L16").Value = 3024.1428
$ws.Range("M16").Value = -2269.2727
$ws.Range("N16").Value = -3364.1428
$ws.Range("H22").Value = 1149.75
$ws.Range("I22").Value = 1149.75
$ws.Range("K22").Value = 1149.75
$ws.Range("M22").Value = -854.75
$ws.Range("H27").Value = 1149.75
$ws.Range("I27").Value = 1149.75
$ws.Range("K27").Value = 1149.75
$ws.Range("M27").Value = -1042.75
$ws.Range("H40").Value = 5997.3335
$ws.Range("I40").Value = 5997.3335
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 5997.3335
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -5861.3335
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 2549.2666
$ws.Range("I46").Value = 2425.818
$ws.Range("J46").Value = 2620.7368
$ws.Range("K46").Value = 2425.818
$ws.Range("L46").Value = 2620.7368
$ws.Range("M46").Value = -2237.818
$ws.Range("N46").Value = -2996.7368
$ws.Range("H55").Value = 83333650
$ws.Range("I55").Value = 83333650
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 83333650
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -83333477
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 4526.143
$ws.Range("I61").Value = 1505.6364
$ws.Range("J61").Value = 15601.333
$ws.Range("K61").Value = 1505.6364
$ws.Range("L61").Value = 15601.333
$ws.Range("M61").Value = -1303.6364
$ws.Range("N61").Value = -16005.333
$ws.Range("H68").Value = 1666.5
$ws.Range("I68").Value = 1399.75
$ws.Range("K68").Value = 1399.75
$ws.Range("M68").Value = -650.75
$ws.Range("H71").Value = 1666.5
$ws.Range("I71").Value = 1399.75
$ws.Range("K71").Value = 6998.75
$ws.Range("M71").Value = -3254.75
$ws.Range("H82").Value = 2510.2
$ws.Range("J82").Value = 2387.25
$ws.Range("L82").Value = 2387.25
$ws.Range("N82").Value = -3109.25
$ws.Range("H85").Value = 2510.2
$ws.Range("J85").Value = 2387.25
$ws.Range("L85").Value = 2387.25
$ws.Range("N85").Value = -4883.25
$ws.Range("H93").Value = 41667704
$ws.Range("I93").Value = 47620044
$ws.Range("K93").Value = 47620044
$ws.Range("M93").Value = -47618796
$ws.Range("H100").Value = 3050.389
$ws.Range("I100").Value = 3167.1333
$ws.Range("K100").Value = 3167.1333
$ws.Range("M100").Value = -2626.1333
$ws.Range("H113").Value = 4526.143
$ws.Range("I113").Value = 1505.6364
$ws.Range("J113").Value = 15601.333
$ws.Range("K113").Value = 1505.6364
$ws.Range("L113").Value = 15601.333
$ws.Range("M113").Value = 664.3635999999999
$ws.Range("N113").Value = -19941.333
$ws.Range("H122").Value = 5617.154
$ws.Range("I122").Value = 4725.5
$ws.Range("K122").Value = 14176.5
$ws.Range("M122").Value = -11726.5
$ws.Range("H126").Value = 63311.35
$ws.Range("I126").Value = 3564.4443
$ws.Range("K126").Value = 10693.3329
$ws.Range("M126").Value = -8223.332900000001
$ws.Range("H132").Value = 127638.44
$ws.Range("I132").Value = 74157.86
$ws.Range("J132").Value = 502002.5
$ws.Range("K132").Value = 222473.58
$ws.Range("L132").Value = 1506007.5
$ws.Range("M132").Value = -219943.58
$ws.Range("N132").Value = -1511067.5
$ws.Range("H136").Value = 183467.58
$ws.Range("I136").Value = 337202.34
$ws.Range("J136").Value = 132222.67
$ws.Range("K136").Value = 1011607.02
$ws.Range("L136").Value = 396668.01
$ws.Range("M136").Value = -1009057.02
$ws.Range("N136").Value = -401768.01

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H74").Value = 25755.834
$ws.Range("J74").Value = 25755.834
$ws.Range("L74").Value = 25755.834
$ws.Range("N74").Value = -27627.834
$ws.Range("H75").Value = 9415634
$ws.Range("J75").Value = 9415634
$ws.Range("L75").Value = 9415634
$ws.Range("N75").Value = -9417506
$ws.Range("H77").Value = 25755.834
$ws.Range("J77").Value = 25755.834
$ws.Range("L77").Value = 77267.50199999999
$ws.Range("N77").Value = -86627.50199999999
$ws.Range("H78").Value = 9415634
$ws.Range("J78").Value = 9415634
$ws.Range("L78").Value = 28246902
$ws.Range("N78").Value = -28256262
$ws.Range("H81").Value = 100500.5
$ws.Range("J81").Value = 1002
$ws.Range("L81").Value = 2004
$ws.Range("N81").Value = -4126
$ws.Range("H84").Value = 100500.5
$ws.Range("J84").Value = 1002
$ws.Range("L84").Value = 10020
$ws.Range("N84").Value = -20628
$ws.Range("H100").Value = 1869.1833
$ws.Range("I100").Value = 1835.2174
$ws.Range("K100").Value = 3670.4348
$ws.Range("M100").Value = -3129.4348
$ws.Range("H107").Value = 21739970
$ws.Range("I107").Value = 38462680
$ws.Range("K107").Value = 115388040
$ws.Range("M107").Value = -115386120
$ws.Range("H113").Value = 759.7143
$ws.Range("I113").Value = 726.3333
$ws.Range("K113").Value = 2178.9999
$ws.Range("M113").Value = -8.999899999999798
$ws.Range("H126").Value = 7756.8623
$ws.Range("I126").Value = 7137.087
$ws.Range("J126").Value = 10132.667
$ws.Range("K126").Value = 21411.261
$ws.Range("L126").Value = 30398.001
$ws.Range("M126").Value = -18941.261
$ws.Range("N126").Value = -35338.001
$ws.Range("H132").Value = 3548.524
$ws.Range("I132").Value = 2854.5293
$ws.Range("K132").Value = 8563.5879
$ws.Range("M132").Value = -6033.5879
$ws.Range("H136").Value = 1694.4546
$ws.Range("I136").Value = 1694.4546
$ws.Range("K136").Value = 5083.3638
$ws.Range("M136").Value = -2533.3638

Write-Host "All changes applied."